$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style/format from G1 (existing "sum" header) to H1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set the new header text and the data value
$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 1
